# mys and sabah data import script updated.
# Append 4 new daily rows (249-252) to the "Data active cases" sheet,
# continuing the existing data series (previously the last data row was 248,
# with rows 249-252 present only as blank, pre-formatted rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Carry the row 248 formatting down onto rows 249:252 for every
#        column, same as dragging the fill handle / copying the row down,
#        so the new rows keep the workbook's existing look (date format on
#        A, centered number styling on B:D/F:G, centered formula styling on
#        E/H, etc.) -----------------------------------------------------
$cols = @("A","B","C","D","E","F","G","H")
foreach ($col in $cols) {
    $ws.Range($col + "248").Copy()
    $ws.Range($col + "249:" + $col + "252").PasteSpecial(-4122)
}

# --- 2. Write the new data values -------------------------------------
# Row 249 : 2020-10-22
$ws.Range("A249").Value = 44126
$ws.Range("B249").Value = 847
$ws.Range("C249").Value = 5
$ws.Range("D249").Value = 8183
$ws.Range("F249").Value = 90
$ws.Range("G249").Value = 29

# Row 250 : 2020-10-23
$ws.Range("A250").Value = 44127
$ws.Range("B250").Value = 710
$ws.Range("C250").Value = 1
$ws.Range("D250").Value = 8416
$ws.Range("F250").Value = 90
$ws.Range("G250").Value = 28

# Row 251 : 2020-10-24
$ws.Range("A251").Value = 44128
$ws.Range("B251").Value = 1228
$ws.Range("C251").Value = 11
$ws.Range("D251").Value = 8966
$ws.Range("F251").Value = 92
$ws.Range("G251").Value = 31

# Row 252 : 2020-10-25
$ws.Range("A252").Value = 44129
$ws.Range("B252").Value = 823
$ws.Range("C252").Value = 0
$ws.Range("D252").Value = 9202
$ws.Range("F252").Value = 99
$ws.Range("G252").Value = 30

# --- 3. Extend the two shared formula columns (E = D - F, H = F / D) down
#        through the new rows, matching how the sheet computes these
#        columns for every other data row --------------------------------
foreach ($r in 249..252) {
    $ws.Range("E" + $r).Formula = "=D" + $r + "-F" + $r
    $ws.Range("H" + $r).Formula = "=F" + $r + "/D" + $r
}

Write-Host "Rows 249-252 populated."
